$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = "If occur license error"
$ws.Range("C21").Value = "run cmd as administrator"
$ws.Range("C22").Value = "cd C:\Users\<user>\AppData\Local\Android\sdk\tools\bin"
$ws.Range("C23").Value = "sdkmanager --licenses"
$ws.Range("C24").Value = "press 'y' until finish"

$ws.Range("C20").Font.Color = 255
